$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2978.8333
$ws.Range("I40").Value = 7933.3335
$ws.Range("J40").Value = 1987.9333
$ws.Range("K40").Value = 7933.3335
$ws.Range("L40").Value = 1987.9333
$ws.Range("M40").Value = -7758.3335
$ws.Range("N40").Value = -2337.9333

$ws.Range("H103").Value = 5437.75
$ws.Range("J103").Value = 5500.5
$ws.Range("L103").Value = 16501.5
$ws.Range("N103").Value = -17673.5

$ws.Range("H112").Value = 1035.8182
$ws.Range("J112").Value = 1061.3334
$ws.Range("L112").Value = 3184.0002
$ws.Range("N112").Value = -5400.0002

$ws.Range("H138").Value = 2786.7656
$ws.Range("I138").Value = 1663.421
$ws.Range("J138").Value = 3261.0667
$ws.Range("K138").Value = 4990.263
$ws.Range("L138").Value = 9783.2001
$ws.Range("M138").Value = 149.7370000000001
$ws.Range("N138").Value = -20063.2001

$ws.Range("H141").Value = 2893.8572
$ws.Range("I141").Value = 2726.1667
$ws.Range("K141").Value = 8178.500100000001
$ws.Range("M141").Value = -2998.500100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3657.71
$ws.Range("I32").Value = 3306.9175
$ws.Range("K32").Value = 3306.9175
$ws.Range("M32").Value = -3019.9175

$ws.Range("H74").Value = 823.439
$ws.Range("I74").Value = 760.0571
$ws.Range("K74").Value = 760.0571
$ws.Range("M74").Value = 113.9429

$ws.Range("H77").Value = 823.439
$ws.Range("I77").Value = 760.0571
$ws.Range("K77").Value = 3800.2855
$ws.Range("M77").Value = 567.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2310.7083
$ws.Range("I80").Value = 901.5
$ws.Range("J80").Value = 3015.3125
$ws.Range("K80").Value = 901.5
$ws.Range("L80").Value = 3015.3125
$ws.Range("M80").Value = 96.5
$ws.Range("N80").Value = -5011.3125

$ws.Range("H83").Value = 2310.7083
$ws.Range("I83").Value = 901.5
$ws.Range("J83").Value = 3015.3125
$ws.Range("K83").Value = 4507.5
$ws.Range("L83").Value = 15076.5625
$ws.Range("M83").Value = 484.5
$ws.Range("N83").Value = -25060.5625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2338.0454
$ws.Range("I58").Value = 2031.4615
$ws.Range("J58").Value = 2780.889
$ws.Range("K58").Value = 2031.4615
$ws.Range("L58").Value = 2780.889
$ws.Range("M58").Value = -1828.4615
$ws.Range("N58").Value = -3186.889

$ws.Range("H136").Value = 2338.0454
$ws.Range("I136").Value = 2031.4615
$ws.Range("J136").Value = 2780.889
$ws.Range("K136").Value = 6094.3845
$ws.Range("L136").Value = 8342.667000000001
$ws.Range("M136").Value = -3544.3845
$ws.Range("N136").Value = -13442.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1199.8667
$ws.Range("J34").Value = 1271.2858
$ws.Range("L34").Value = 3813.8574
$ws.Range("N34").Value = -3981.8574

$ws.Range("H39").Value = 2148.5
$ws.Range("J39").Value = 2698
$ws.Range("L39").Value = 8094
$ws.Range("N39").Value = -8682

$ws.Range("H55").Value = 17221.111
$ws.Range("I55").Value = 25572.5
$ws.Range("J55").Value = 10540
$ws.Range("K55").Value = 76717.5
$ws.Range("L55").Value = 31620
$ws.Range("M55").Value = -76540.5
$ws.Range("N55").Value = -31974

$ws.Range("H113").Value = 1069.5264
$ws.Range("I113").Value = 2133.8333
$ws.Range("J113").Value = 578.3077
$ws.Range("K113").Value = 6401.499899999999
$ws.Range("L113").Value = 1734.9231
$ws.Range("M113").Value = -4231.499899999999
$ws.Range("N113").Value = -6074.9231

$ws.Range("H131").Value = 854.04
$ws.Range("J131").Value = 858.2041
$ws.Range("L131").Value = 2574.6123
$ws.Range("N131").Value = -12654.6123

$ws.Range("H134").Value = 4003.5
$ws.Range("I134").Value = 2536.3635
$ws.Range("J134").Value = 4852.8945
$ws.Range("K134").Value = 7609.0905
$ws.Range("L134").Value = 14558.6835
$ws.Range("M134").Value = -2539.0905
$ws.Range("N134").Value = -24698.6835

$ws.Range("H137").Value = 50344.87
$ws.Range("I137").Value = 60525.824
$ws.Range("J137").Value = 21498.834
$ws.Range("K137").Value = 181577.472
$ws.Range("L137").Value = 64496.50199999999
$ws.Range("M137").Value = -176477.472
$ws.Range("N137").Value = -74696.50199999999

$ws.Range("H138").Value = 8032
$ws.Range("I138").Value = 12746.667
$ws.Range("J138").Value = 1970.2858
$ws.Range("K138").Value = 38240.001
$ws.Range("L138").Value = 5910.857400000001
$ws.Range("M138").Value = -33100.001
$ws.Range("N138").Value = -16190.8574

$ws.Range("H139").Value = 2029.4642
$ws.Range("I139").Value = 1204.375
$ws.Range("J139").Value = 3129.5833
$ws.Range("K139").Value = 3613.125
$ws.Range("L139").Value = 9388.749899999999
$ws.Range("M139").Value = 1526.875
$ws.Range("N139").Value = -19668.7499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1206.8572
$ws.Range("I122").Value = 1446.6666
$ws.Range("K122").Value = 4339.9998
$ws.Range("M122").Value = -1889.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 675919.9399999999
$ws.Range("I46").Value = 343.33334
$ws.Range("J46").Value = 1689284.9
$ws.Range("K46").Value = 343.33334
$ws.Range("L46").Value = 1689284.9
$ws.Range("M46").Value = -155.33334
$ws.Range("N46").Value = -1689660.9

$ws.Range("H55").Value = 292086.78
$ws.Range("I55").Value = 517312.1
$ws.Range("J55").Value = 618.7059
$ws.Range("K55").Value = 517312.1
$ws.Range("L55").Value = 618.7059
$ws.Range("M55").Value = -517139.1
$ws.Range("N55").Value = -964.7059

$ws.Range("H68").Value = 4600
$ws.Range("I68").Value = 1700
$ws.Range("J68").Value = 5014.2856
$ws.Range("K68").Value = 1700
$ws.Range("L68").Value = 5014.2856
$ws.Range("M68").Value = -951
$ws.Range("N68").Value = -6512.2856

$ws.Range("H71").Value = 4600
$ws.Range("I71").Value = 1700
$ws.Range("J71").Value = 5014.2856
$ws.Range("K71").Value = 8500
$ws.Range("L71").Value = 25071.428
$ws.Range("M71").Value = -4756
$ws.Range("N71").Value = -32559.428

$ws.Range("H132").Value = 5260.7144
$ws.Range("I132").Value = 5007.25
$ws.Range("J132").Value = 5598.6665
$ws.Range("K132").Value = 15021.75
$ws.Range("L132").Value = 16795.9995
$ws.Range("M132").Value = -12491.75
$ws.Range("N132").Value = -21855.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6995259.5
$ws.Range("I62").Value = 12822392
$ws.Range("K62").Value = 12822392
$ws.Range("M62").Value = -12821768

$ws.Range("H65").Value = 6995259.5
$ws.Range("I65").Value = 12822392
$ws.Range("K65").Value = 64111960
$ws.Range("M65").Value = -64108840

$ws.Range("H122").Value = 2744.0557
$ws.Range("I122").Value = 1948.5834
$ws.Range("J122").Value = 4335
$ws.Range("K122").Value = 5845.7502
$ws.Range("L122").Value = 13005
$ws.Range("M122").Value = -3395.7502
$ws.Range("N122").Value = -17905

$ws.Range("H136").Value = 1785.9215
$ws.Range("I136").Value = 746.5714
$ws.Range("J136").Value = 2513.4666
$ws.Range("K136").Value = 2239.7142
$ws.Range("L136").Value = 7540.399800000001
$ws.Range("M136").Value = 310.2857999999997
$ws.Range("N136").Value = -12640.3998
